# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" sheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 554
    "F3" = 3544
    "F4" = 99
    "F5" = 688
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
